# Update "想去人数" (F) counts and a couple of "最低票价" (G) sold-out
# markers across the 展览 / 演出 / 全部类型 sheets to match the refreshed
# scrape output ("Update gh-pages to output generated at 456a3b4").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 52
$ws.Range("F4").Value = 3062
$ws.Range("F5").Value = 1688
$ws.Range("F6").Value = 2068
$ws.Range("F7").Value = 324
$ws.Range("F8").Value = 303
$ws.Range("F9").Value = 884
$ws.Range("F10").Value = 970
$ws.Range("F11").Value = 213
$ws.Range("F12").Value = 437
$ws.Range("F13").Value = 1140
$ws.Range("F15").Value = 71
$ws.Range("F17").Value = 7415
$ws.Range("F18").Value = 305
$ws.Range("F19").Value = 2447
$ws.Range("G19").Value = "已售罄"
$ws.Range("F20").Value = 198
$ws.Range("F21").Value = 212
$ws.Range("F22").Value = 166
$ws.Range("F23").Value = 450
$ws.Range("F24").Value = 514
$ws.Range("F25").Value = 73
$ws.Range("F26").Value = 1127
$ws.Range("F27").Value = 956
$ws.Range("F29").Value = 344
$ws.Range("F30").Value = 240
$ws.Range("F31").Value = 1131
$ws.Range("F32").Value = 1907
$ws.Range("F34").Value = 25
$ws.Range("F35").Value = 157
$ws.Range("F37").Value = 35
$ws.Range("F38").Value = 155
$ws.Range("F39").Value = 301
$ws.Range("F41").Value = 198
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 19
$ws.Range("F6").Value = 10
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 19
$ws.Range("F6").Value = 52
$ws.Range("F7").Value = 3062
$ws.Range("F8").Value = 1688
$ws.Range("F9").Value = 2068
$ws.Range("F10").Value = 324
$ws.Range("F11").Value = 303
$ws.Range("F12").Value = 884
$ws.Range("F14").Value = 970
$ws.Range("F15").Value = 213
$ws.Range("F16").Value = 437
$ws.Range("F17").Value = 1140
$ws.Range("F19").Value = 71
$ws.Range("F21").Value = 7417
$ws.Range("F22").Value = 305
$ws.Range("F23").Value = 2447
$ws.Range("G23").Value = "已售罄"
$ws.Range("F24").Value = 10
$ws.Range("F25").Value = 198
$ws.Range("F26").Value = 212
$ws.Range("F27").Value = 166
$ws.Range("F28").Value = 450
$ws.Range("F29").Value = 514
$ws.Range("F30").Value = 73
$ws.Range("F31").Value = 1127
$ws.Range("F32").Value = 956
$ws.Range("F34").Value = 345
$ws.Range("F35").Value = 240
$ws.Range("F36").Value = 1131
$ws.Range("F37").Value = 1907
$ws.Range("F39").Value = 25
$ws.Range("F40").Value = 157
$ws.Range("F42").Value = 35
$ws.Range("F43").Value = 155
$ws.Range("F44").Value = 301
$ws.Range("F49").Value = 198
